$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-10-04 18:21:28"

$ws.Range("A2").Value = $newTimestamp
$ws.Range("A3").Value = $newTimestamp
$ws.Range("A4").Value = $newTimestamp
$ws.Range("A5").Value = $newTimestamp
$ws.Range("A6").Value = $newTimestamp
